# "Crea fichero elementos completo"
# Inserts four new parameter rows (CODIGO_DELEGACION, IDIOMA_EXPEDIENTE,
# MEDIO_NOTIFICACION, INTERESADO_NOTIFICACION) above the existing
# "*Director" rows, widens column A, and moves the active selection —
# matching the author's final edit of CT102A_datosEspecificos.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows (old rows 2-5) down to rows 6-9, opening up
# four blank rows right under the header row.
$ws.Rows.Item(2).Resize(4).Insert()

# The newly-inserted rows pick up the bold header formatting from row 1
# (Excel's default "format from above" behaviour on row insert) and are
# fully populated with blank cells across A:F. Strip that back down: no
# bold style, and only the "key" (A) / "path" (E) columns actually hold
# data, same as the rest of the sheet's new rows.
$ws.Range("A2:F5").ClearFormats()
$ws.Range("B2:D5").ClearContents()
$ws.Range("F2:F5").ClearContents()

# Fill in column A ("key") first for all four new rows ...
$ws.Range("A2").Value = "CODIGO_DELEGACION"
$ws.Range("A3").Value = "IDIOMA_EXPEDIENTE"
$ws.Range("A4").Value = "MEDIO_NOTIFICACION"
$ws.Range("A5").Value = "INTERESADO_NOTIFICACION"

# ... then column E ("path") for all four, so shared strings land in the
# same order the author's workbook used.
$ws.Range("E2").Value = "//ProcedimientoXunta/SI460A_4/Delegacion/cmbDelegacion"
$ws.Range("E3").Value = "//ProcedimientoXunta/SI460A_4/Delegacion/txtIdioma"
$ws.Range("E4").Value = "//ProcedimientoXunta/SI460A/Notificacion/rblModalidad"
$ws.Range("E5").Value = "//ProcedimientoXunta/SI460A/Notificacion/rbNotificar"

# Widen column A to fit the new, longer keys, and leave the selection where
# the author left it when they saved the file.
$ws.Columns.Item(1).ColumnWidth = 27.3
[void]$ws.Range("B16").Select()
